# Applies the cryptocurrency price/volume refresh described by the commit diff.
# Cells D2:E51 are free-form text (prices like "67.134.26" / "2.467.96" are
# formatted strings, not real numbers) so any new value that *looks* purely
# numeric is pinned to the Text number format before being written - otherwise
# Excel would silently reinterpret it (and e.g. "1.00" would collapse to "1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.121.30'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '2.467.84'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '582.57'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '174.51'
$ws.Range('E6').Value = '  +3.23%  '
$ws.Range('E8').Value = '  -0.48%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.137'
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('E10').Value = '  +0.46%  '
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('E14').Value = '  -0.94%  '
$ws.Range('D15').Value = '67.021.91'
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('E16').Value = '  -0.39%  '
$ws.Range('D17').Value = '2.421.28'
$ws.Range('E17').Value = '  -2.41%  '
$ws.Range('E18').Value = '  -2.66%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.45'
$ws.Range('E19').Value = '  -2.05%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '348.94'
$ws.Range('E20').Value = '  -1.49%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '3.98'
$ws.Range('E21').Value = '  -1.47%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '69.28'
$ws.Range('E23').Value = '  +0.26%  '
$ws.Range('E24').Value = '  -1.16%  '
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.12'
$ws.Range('E26').Value = '  -1.98%  '
$ws.Range('D27').Value = '2.594.29'
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('D29').Value = '0.0₃0899'
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '499.32'
$ws.Range('E30').Value = '  -3.60%  '
$ws.Range('E31').Value = '  -0.46%  '
$ws.Range('E32').Value = '  -1.15%  '
$ws.Range('E33').Value = '  -1.70%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '161.99'
$ws.Range('E36').Value = '  +2.56%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '18.67'
$ws.Range('E38').Value = '  -1.63%  '
$ws.Range('E39').Value = '  -2.21%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  +1.39%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.326'
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.82'
$ws.Range('E43').Value = '  +0.41%  '
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '142.01'
$ws.Range('E45').Value = '  +0.55%  '
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('E47').Value = '  -0.92%  '
$ws.Range('D48').Value = '0.0₆0252'
$ws.Range('E48').Value = '  -1.23%  '
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('E50').Value = '  -1.64%  '
$ws.Range('E51').Value = '  -0.07%  '
